$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2941.22
$ws.Cells.Item(17, 10).Value = 2941.22
$ws.Cells.Item(17, 12).Value = 8823.66
$ws.Cells.Item(17, 14).Value = -9159.66
$ws.Cells.Item(92, 8).Value = 468.5
$ws.Cells.Item(92, 9).Value = 460.57144
$ws.Cells.Item(92, 10).Value = 496.25
$ws.Cells.Item(92, 11).Value = 460.57144
$ws.Cells.Item(92, 12).Value = 496.25
$ws.Cells.Item(92, 13).Value = 787.4285600000001
$ws.Cells.Item(92, 14).Value = -2992.25
$ws.Cells.Item(96, 8).Value = 41671604
$ws.Cells.Item(96, 9).Value = 2917.077
$ws.Cells.Item(96, 10).Value = 90916420
$ws.Cells.Item(96, 11).Value = 8751.231
$ws.Cells.Item(96, 12).Value = 272749260
$ws.Cells.Item(96, 13).Value = -7378.231
$ws.Cells.Item(96, 14).Value = -272752006
$ws.Cells.Item(137, 8).Value = 3321.9138
$ws.Cells.Item(137, 9).Value = 1082.2727
$ws.Cells.Item(137, 10).Value = 10360.786
$ws.Cells.Item(137, 11).Value = 3246.8181
$ws.Cells.Item(137, 12).Value = 31082.358
$ws.Cells.Item(137, 13).Value = -696.8181
$ws.Cells.Item(137, 14).Value = -36182.358
$ws.Cells.Item(138, 8).Value = 1863.2
$ws.Cells.Item(138, 9).Value = 602.7273
$ws.Cells.Item(138, 10).Value = 2853.5715
$ws.Cells.Item(138, 11).Value = 1808.1819
$ws.Cells.Item(138, 12).Value = 8560.7145
$ws.Cells.Item(138, 13).Value = 3331.8181
$ws.Cells.Item(138, 14).Value = -18840.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2041.8
$ws.Cells.Item(45, 9).Value = 1636.3478
$ws.Cells.Item(45, 10).Value = 2818.9167
$ws.Cells.Item(45, 11).Value = 1636.3478
$ws.Cells.Item(45, 12).Value = 2818.9167
$ws.Cells.Item(45, 13).Value = -1259.3478
$ws.Cells.Item(45, 14).Value = -3572.9167
$ws.Cells.Item(61, 8).Value = 1378.0615
$ws.Cells.Item(61, 9).Value = 1217.54
$ws.Cells.Item(61, 10).Value = 1913.1333
$ws.Cells.Item(61, 11).Value = 1217.54
$ws.Cells.Item(61, 12).Value = 1913.1333
$ws.Cells.Item(61, 13).Value = -1005.54
$ws.Cells.Item(61, 14).Value = -2337.1333
$ws.Cells.Item(74, 8).Value = 1716.3392
$ws.Cells.Item(74, 9).Value = 1734.1464
$ws.Cells.Item(74, 11).Value = 1734.1464
$ws.Cells.Item(74, 13).Value = -860.1464000000001
$ws.Cells.Item(77, 8).Value = 1716.3392
$ws.Cells.Item(77, 9).Value = 1734.1464
$ws.Cells.Item(77, 11).Value = 8670.732
$ws.Cells.Item(77, 13).Value = -4302.732
$ws.Cells.Item(102, 8).Value = 10863.462
$ws.Cells.Item(102, 9).Value = 2179.4
$ws.Cells.Item(102, 11).Value = 2179.4
$ws.Cells.Item(102, 13).Value = -557.4000000000001
$ws.Cells.Item(110, 8).Value = 1477.625
$ws.Cells.Item(110, 9).Value = 1481.2778
$ws.Cells.Item(110, 10).Value = 1466.6666
$ws.Cells.Item(110, 11).Value = 1481.2778
$ws.Cells.Item(110, 12).Value = 1466.6666
$ws.Cells.Item(110, 13).Value = 563.7221999999999
$ws.Cells.Item(110, 14).Value = -5556.6666
$ws.Cells.Item(132, 8).Value = 12197355
$ws.Cells.Item(132, 9).Value = 18519966
$ws.Cells.Item(132, 10).Value = 3744.9285
$ws.Cells.Item(132, 11).Value = 55559898
$ws.Cells.Item(132, 12).Value = 11234.7855
$ws.Cells.Item(132, 13).Value = -55557368
$ws.Cells.Item(132, 14).Value = -16294.7855
$ws.Cells.Item(136, 8).Value = 1378.0615
$ws.Cells.Item(136, 9).Value = 1217.54
$ws.Cells.Item(136, 10).Value = 1913.1333
$ws.Cells.Item(136, 11).Value = 3652.62
$ws.Cells.Item(136, 12).Value = 5739.3999
$ws.Cells.Item(136, 13).Value = -1102.62
$ws.Cells.Item(136, 14).Value = -10839.3999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2155.077
$ws.Cells.Item(99, 9).Value = 2155.077
$ws.Cells.Item(99, 11).Value = 2155.077
$ws.Cells.Item(99, 13).Value = -657.0770000000002
$ws.Cells.Item(105, 8).Value = 2171.353
$ws.Cells.Item(105, 9).Value = 2730.8
$ws.Cells.Item(105, 10).Value = 2074.8965
$ws.Cells.Item(105, 11).Value = 2730.8
$ws.Cells.Item(105, 12).Value = 2074.8965
$ws.Cells.Item(105, 13).Value = -983.8000000000002
$ws.Cells.Item(105, 14).Value = -5568.8965
$ws.Cells.Item(107, 8).Value = 2524.4707
$ws.Cells.Item(107, 9).Value = 2355.5454
$ws.Cells.Item(107, 10).Value = 2834.1667
$ws.Cells.Item(107, 11).Value = 2355.5454
$ws.Cells.Item(107, 12).Value = 2834.1667
$ws.Cells.Item(107, 13).Value = -435.5454
$ws.Cells.Item(107, 14).Value = -6674.1667
$ws.Cells.Item(134, 8).Value = 209880.4
$ws.Cells.Item(134, 9).Value = 2745.6667
$ws.Cells.Item(134, 10).Value = 241477.23
$ws.Cells.Item(134, 11).Value = 8237.000100000001
$ws.Cells.Item(134, 12).Value = 724431.6900000001
$ws.Cells.Item(134, 13).Value = -5702.000100000001
$ws.Cells.Item(134, 14).Value = -729501.6900000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1397.8
$ws.Cells.Item(16, 9).Value = 1330
$ws.Cells.Item(16, 10).Value = 1499.5
$ws.Cells.Item(16, 11).Value = 1330
$ws.Cells.Item(16, 12).Value = 1499.5
$ws.Cells.Item(16, 13).Value = -1043
$ws.Cells.Item(16, 14).Value = -2073.5
$ws.Cells.Item(31, 8).Value = 2125.43
$ws.Cells.Item(31, 10).Value = 3266.5305
$ws.Cells.Item(31, 12).Value = 3266.5305
$ws.Cells.Item(31, 14).Value = -3856.5305
$ws.Cells.Item(34, 8).Value = 2125.43
$ws.Cells.Item(34, 10).Value = 3266.5305
$ws.Cells.Item(34, 12).Value = 3266.5305
$ws.Cells.Item(34, 14).Value = -3670.5305
$ws.Cells.Item(58, 8).Value = 1214.9038
$ws.Cells.Item(58, 9).Value = 836.9524
$ws.Cells.Item(58, 11).Value = 836.9524
$ws.Cells.Item(58, 13).Value = -633.9524
$ws.Cells.Item(99, 8).Value = 2154.125
$ws.Cells.Item(99, 9).Value = 2131
$ws.Cells.Item(99, 10).Value = 2177.25
$ws.Cells.Item(99, 11).Value = 2131
$ws.Cells.Item(99, 12).Value = 2177.25
$ws.Cells.Item(99, 13).Value = -633
$ws.Cells.Item(99, 14).Value = -5173.25
$ws.Cells.Item(113, 8).Value = 1397.8
$ws.Cells.Item(113, 9).Value = 1330
$ws.Cells.Item(113, 10).Value = 1499.5
$ws.Cells.Item(113, 11).Value = 1330
$ws.Cells.Item(113, 12).Value = 1499.5
$ws.Cells.Item(113, 13).Value = 840
$ws.Cells.Item(113, 14).Value = -5839.5
$ws.Cells.Item(126, 8).Value = 2154.125
$ws.Cells.Item(126, 9).Value = 2131
$ws.Cells.Item(126, 10).Value = 2177.25
$ws.Cells.Item(126, 11).Value = 6393
$ws.Cells.Item(126, 12).Value = 6531.75
$ws.Cells.Item(126, 13).Value = -3923
$ws.Cells.Item(126, 14).Value = -11471.75
$ws.Cells.Item(132, 8).Value = 46158.03
$ws.Cells.Item(132, 9).Value = 2497.35
$ws.Cells.Item(132, 10).Value = 118925.836
$ws.Cells.Item(132, 11).Value = 7492.049999999999
$ws.Cells.Item(132, 12).Value = 356777.508
$ws.Cells.Item(132, 13).Value = -4962.049999999999
$ws.Cells.Item(132, 14).Value = -361837.508
$ws.Cells.Item(134, 8).Value = 468612.03
$ws.Cells.Item(134, 9).Value = 1144.8889
$ws.Cells.Item(134, 10).Value = 1169812.8
$ws.Cells.Item(134, 11).Value = 3434.6667
$ws.Cells.Item(134, 12).Value = 3509438.4
$ws.Cells.Item(134, 13).Value = -899.6666999999998
$ws.Cells.Item(134, 14).Value = -3514508.4
$ws.Cells.Item(136, 8).Value = 1214.9038
$ws.Cells.Item(136, 9).Value = 836.9524
$ws.Cells.Item(136, 11).Value = 2510.8572
$ws.Cells.Item(136, 13).Value = 39.14280000000008

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2982.6
$ws.Cells.Item(5, 9).Value = 5305.8096
$ws.Cells.Item(5, 10).Value = 1300.2759
$ws.Cells.Item(5, 11).Value = 15917.4288
$ws.Cells.Item(5, 12).Value = 3900.8277
$ws.Cells.Item(5, 13).Value = -15805.4288
$ws.Cells.Item(5, 14).Value = -4124.8277
$ws.Cells.Item(98, 8).Value = 1040.7142
$ws.Cells.Item(98, 9).Value = 495.33334
$ws.Cells.Item(98, 10).Value = 1449.75
$ws.Cells.Item(98, 11).Value = 1486.00002
$ws.Cells.Item(98, 12).Value = 4349.25
$ws.Cells.Item(98, 13).Value = 11.99998000000005
$ws.Cells.Item(98, 14).Value = -7345.25
$ws.Cells.Item(122, 8).Value = 2046.3522
$ws.Cells.Item(122, 9).Value = 586.34045
$ws.Cells.Item(122, 10).Value = 4905.5415
$ws.Cells.Item(122, 11).Value = 5277.06405
$ws.Cells.Item(122, 12).Value = 44149.8735
$ws.Cells.Item(122, 13).Value = -2827.06405
$ws.Cells.Item(122, 14).Value = -49049.8735
$ws.Cells.Item(135, 8).Value = 2982.6
$ws.Cells.Item(135, 9).Value = 5305.8096
$ws.Cells.Item(135, 10).Value = 1300.2759
$ws.Cells.Item(135, 11).Value = 47752.2864
$ws.Cells.Item(135, 12).Value = 11702.4831
$ws.Cells.Item(135, 13).Value = -45217.2864
$ws.Cells.Item(135, 14).Value = -16772.4831
$ws.Cells.Item(137, 8).Value = 38470292
$ws.Cells.Item(137, 9).Value = 3044.6
$ws.Cells.Item(137, 10).Value = 90925630
$ws.Cells.Item(137, 11).Value = 9133.799999999999
$ws.Cells.Item(137, 12).Value = 272776890
$ws.Cells.Item(137, 13).Value = -4033.799999999999
$ws.Cells.Item(137, 14).Value = -272787090

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 5429.727
$ws.Cells.Item(97, 9).Value = 3088.125
$ws.Cells.Item(97, 10).Value = 11674
$ws.Cells.Item(97, 11).Value = 3088.125
$ws.Cells.Item(97, 12).Value = 11674
$ws.Cells.Item(97, 13).Value = -2592.125
$ws.Cells.Item(97, 14).Value = -12666
$ws.Cells.Item(102, 8).Value = 1423.68
$ws.Cells.Item(102, 9).Value = 1454.2273
$ws.Cells.Item(102, 10).Value = 1199.6666
$ws.Cells.Item(102, 11).Value = 1454.2273
$ws.Cells.Item(102, 12).Value = 1199.6666
$ws.Cells.Item(102, 13).Value = 167.7727
$ws.Cells.Item(102, 14).Value = -4443.6666
$ws.Cells.Item(126, 8).Value = 6318.3477
$ws.Cells.Item(126, 9).Value = 9859.691999999999
$ws.Cells.Item(126, 10).Value = 1714.6
$ws.Cells.Item(126, 11).Value = 29579.076
$ws.Cells.Item(126, 12).Value = 5143.799999999999
$ws.Cells.Item(126, 13).Value = -27109.076
$ws.Cells.Item(126, 14).Value = -10083.8
$ws.Cells.Item(132, 8).Value = 2681.2856
$ws.Cells.Item(132, 9).Value = 1856.85
$ws.Cells.Item(132, 10).Value = 3780.5334
$ws.Cells.Item(132, 11).Value = 5570.549999999999
$ws.Cells.Item(132, 12).Value = 11341.6002
$ws.Cells.Item(132, 13).Value = -3040.549999999999
$ws.Cells.Item(132, 14).Value = -16401.6002

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2525.2744
$ws.Cells.Item(132, 9).Value = 1815.0769
$ws.Cells.Item(132, 10).Value = 4833.4165
$ws.Cells.Item(132, 11).Value = 5445.2307
$ws.Cells.Item(132, 12).Value = 14500.2495
$ws.Cells.Item(132, 13).Value = -2915.2307
$ws.Cells.Item(132, 14).Value = -19560.2495

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(108, 8).Value = 44622
$ws.Cells.Item(108, 10).Value = 44622
$ws.Cells.Item(108, 12).Value = 44622
$ws.Cells.Item(108, 14).Value = -52302
$ws.Cells.Item(132, 8).Value = 2039.262
$ws.Cells.Item(132, 9).Value = 2056.4814
$ws.Cells.Item(132, 10).Value = 2008.2667
$ws.Cells.Item(132, 11).Value = 6169.4442
$ws.Cells.Item(132, 12).Value = 6024.800099999999
$ws.Cells.Item(132, 13).Value = -3639.4442
$ws.Cells.Item(132, 14).Value = -11084.8001

Write-Output "Applied 243 cell edits"